$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bracket")

$updates = @{
    "D2"  = "Four-Striped Grass Mouse"
    "E4"  = "Giant striped mongoose"
    "F8"  = "Striped dolphin"
    "L8"  = "Golden Eagle"
    "D10" = "Side-striped jackal"
    "M12" = "Tent-making Bat"
    "D14" = "Striped dolphin"
    "N14" = "Tent-making Bat"
    "G16" = "Kudu"
    "E20" = "Wildcat"
    "N22" = "Homo habilis"
    "F24" = "Kudu"
    "N26" = "Montezuma Oropendola"
    "H32" = "Rock Hyrax "
    "I32" = "Greater Rhea"
    "D34" = "Sea Otter"
    "M36" = "Caspian Terns"
    "F40" = "Grasshopper Mouse"
    "L40" = "Siamang"
    "D42" = "Sibree Dwarf Lemur"
    "E44" = "Siberian Chipmunk"
    "M44" = "Siamang"
    "N46" = "Siamang"
    "G48" = "Rock Hyrax "
    "E52" = "Dik Dik"
    "M52" = "Giant Water Bug"
    "D54" = "Colo Colo Opossum"
    "N54" = "Giant Water Bug"
    "F56" = "Rock Hyrax "
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
